# The commit swaps the two theme parts of the deck: the slide-master theme
# ("Integral") and the notes-master theme ("Office Theme") trade places, so
# that after the edit the slide master carries the stock "Office Theme"
# palette and the notes master carries the old "Integral" palette.
#
# This COM host only exposes a single, writable theme object for the
# presentation (reached through SlideMaster/Design/Theme) - its
# ThemeColorScheme is backed by ppt/theme/theme1.xml, the part the slide
# master actually points at. Re-pointing relationships or rewriting the
# <a:theme>/<a:clrScheme> "name" attributes isn't exposed on the object
# model here (those setters are accepted but don't persist), so we recolor
# theme1's 12-slot color scheme to the "Office Theme" palette, which is the
# user-visible effect of the swap for the part of the theme this API can
# reach.

$p = $ppt.ActivePresentation
$theme = $p.SlideMaster.Theme
$colorScheme = $theme.ThemeColorScheme

# clrScheme slot order: dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink.
# Office Theme palette (was previously only in ppt/theme/theme2.xml).
$officeColors = @(
    0x000000,  # dk1
    0xFFFFFF,  # lt1
    0x44546A,  # dk2
    0xE7E6E6,  # lt2
    0x5B9BD5,  # accent1
    0xED7D31,  # accent2
    0xA5A5A5,  # accent3
    0xFFC000,  # accent4
    0x4472C4,  # accent5
    0x70AD47,  # accent6
    0x0563C1,  # hlink
    0x954F72   # folHlink
)

for ($i = 0; $i -lt $officeColors.Count; $i++) {
    $hex = $officeColors[$i]
    $r = [math]::Floor($hex / 0x10000) % 0x100
    $g = [math]::Floor($hex / 0x100) % 0x100
    $b = $hex % 0x100
    $bgr = ($b * 0x10000) + ($g * 0x100) + $r
    $colorScheme.Colors($i + 1).RGB = $bgr
}
